$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.797.59'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.51%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.632.90'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.58%  '
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.93'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5023'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.41%  '
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2571'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06405'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.52%  '
$ws.Range('E10').Value = '  -2.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07664'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.654.07'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.234'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.857.92'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5456'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅7923'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.78%  '
$ws.Range('E17').Value = '  -1.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.828.56'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.002'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '202.56'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.304'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.928'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.947'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.53%  '
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.932'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +10.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.07'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1143'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.03%  '
$ws.Range('E28').Value = '  -0.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.684'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.240'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04970'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.274'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.184'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.531'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.51%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.350'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.171.83'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.8916'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.619'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5570'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.28%  '
$ws.Range('E40').Value = '  -2.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.558'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.69%  '
$ws.Range('E42').Value = '  -0.28%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.640'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.61%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.25'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8021'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.770.54'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.43%  '
$ws.Range('E47').Value = '  -0.11%  '
$ws.Range('E48').Value = '  -0.54%  '
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '54.76'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05027'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.71%  '
